$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range('A11').Value = 'Caffe Nero Glasgow Union St'
$ws.Range('B11').Value = 'https://uk.indeed.com/viewjob?jk=3c30a07d77d3a479'
$ws.Range('C11').Value = 'Posted'
$ws.Range('D11').Value = 'As an Assistant Manager you will support in the running of the store, providing great customer service and the best standards of coffee and food.'
$ws.Range('E11').Value = 44721
$ws.Range('E11').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F11').ClearContents()
$ws.Range('G11').ClearContents()
$ws.Range('H11').Value = 'Glasgow'
$ws.Range('I11').Value = 'G1'

# Row 12
$ws.Range('A12').Value = 'BRGR'
$ws.Range('B12').Value = 'https://uk.indeed.com/viewjob?jk=ceb8df5a07ba8a1b'
$ws.Range('C12').Value = 'Posted'
$ws.Range('D12').Value = 'Day to day responsibility and accountability for the venue in the absence of the manager.'
$ws.Range('E12').Value = 44721
$ws.Range('E12').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F12').Value = 24000
$ws.Range('G12').Value = 24000
$ws.Range('H12').Value = 'Glasgow'
$ws.Range('I12').Value = 'G1'

# Row 13
$ws.Range('A13').Value = 'Frankie and Benny''s'
$ws.Range('B13').Value = 'https://uk.indeed.com/viewjob?jk=89956d0b5f22bce6'
$ws.Range('C13').Value = 'Posted'
$ws.Range('D13').Value = 'Refer a friend scheme with generous bonuses for each successful referral.'
$ws.Range('E13').Value = 44790
$ws.Range('E13').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F13').Value = 27000
$ws.Range('G13').Value = 27000
$ws.Range('H13').Value = 'Glasgow'
$ws.Range('I13').Value = 'G33'

# Row 14
$ws.Range('A14').Value = 'Subway'
$ws.Range('B14').Value = 'https://uk.indeed.com/viewjob?jk=0e6095c2b85ef6c5'
$ws.Range('C14').Value = 'Active'
$ws.Range('D14').Value = 'View all Subway jobs - Bathgate jobs'
$ws.Range('E14').Value = 44792
$ws.Range('E14').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F14').ClearContents()
$ws.Range('G14').ClearContents()
$ws.Range('H14').Value = 'Bathgate'
$ws.Range('I14').Value = 'H48'

# Row 15
$ws.Range('A15').Value = 'Greene King'
$ws.Range('B15').Value = 'https://uk.indeed.com/viewjob?jk=9ca24025cd40b6b4'
$ws.Range('C15').Value = 'Posted'
$ws.Range('D15').Value = 'Competitive salary with generous bonus and pension contribution.'
$ws.Range('E15').Value = 44792
$ws.Range('E15').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F15').Value = 25160
$ws.Range('G15').Value = 27000
$ws.Range('H15').Value = 'Glasgow'
$ws.Range('I15').ClearContents()

# Row 16
$ws.Range('A16').Value = 'Domino''s'
$ws.Range('B16').Value = 'https://uk.indeed.com/viewjob?jk=827e4b72fc5c31c5'
$ws.Range('C16').Value = 'Posted'
$ws.Range('D16').Value = 'Assisting and deputising for the Store Manager, you will receive training and development in all areas of store operations, learning about the daily aspects of…'
$ws.Range('E16').Value = 44798
$ws.Range('E16').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F16').Value = 11
$ws.Range('G16').Value = 11
$ws.Range('H16').Value = 'Carluke'
$ws.Range('I16').ClearContents()

# Row 17
$ws.Range('A17').Value = 'Nando’s Chickenland Limited'
$ws.Range('B17').Value = 'https://uk.indeed.com/viewjob?jk=e436861aec7e6027'
$ws.Range('C17').Value = 'Posted'
$ws.Range('D17').Value = 'An Assistant Manager with previous management experience and bucket-loads of passion and energy to inspire, motivate, and engage your team.'
$ws.Range('E17').Value = 44799
$ws.Range('E17').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F17').Value = 27000
$ws.Range('G17').Value = 27000
$ws.Range('H17').Value = 'Falkirk'
$ws.Range('I17').Value = 'K1'

# Row 18
$ws.Range('A18').Value = 'Hospo Ltd'
$ws.Range('B18').Value = 'https://uk.indeed.com/viewjob?jk=61183ed69de1164c'
$ws.Range('C18').Value = 'Posted'
$ws.Range('D18').Value = 'Handling customer complaints or concerns regarding service or drink quality in a courteous manner.'
$ws.Range('E18').Value = 44800
$ws.Range('E18').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F18').Value = 28000
$ws.Range('G18').Value = 30000
$ws.Range('H18').Value = 'Greenock'
$ws.Range('I18').ClearContents()

# Row 19
$ws.Range('A19').Value = 'Hospo Ltd'
$ws.Range('B19').Value = 'https://uk.indeed.com/viewjob?jk=e9e5cc46eba49088'
$ws.Range('C19').Value = 'Posted'
$ws.Range('D19').Value = 'Helping general manager with the day to day.'
$ws.Range('E19').Value = 44800
$ws.Range('E19').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F19').Value = 30000
$ws.Range('G19').Value = 35000
$ws.Range('H19').Value = 'Greenock'
$ws.Range('I19').ClearContents()

# Row 20
$ws.Range('A20').Value = 'C&C Group'
$ws.Range('B20').Value = 'https://uk.indeed.com/viewjob?jk=2c8dd70f8caa8462'
$ws.Range('C20').Value = 'Active'
$ws.Range('D20').Value = 'Greek food experience is preferred but not essential as full training will be provided.'
$ws.Range('E20').Value = 44804
$ws.Range('E20').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F20').Value = 10
$ws.Range('G20').Value = 12
$ws.Range('H20').Value = 'Glasgow'
$ws.Range('I20').Value = 'G12'

# Row 21
$ws.Range('A21').Value = 'Miller & Carter'
$ws.Range('B21').Value = 'https://uk.indeed.com/viewjob?jk=1522e0d407906305'
$ws.Range('C21').Value = 'Active'
$ws.Range('D21').Value = 'You''ll pull everything together to make sure our guests, and teams, have a great time.'
$ws.Range('E21').Value = 44804
$ws.Range('E21').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F21').Value = 30000
$ws.Range('G21').Value = 30000
$ws.Range('H21').Value = 'Newton'
$ws.Range('I21').ClearContents()

# Row 22
$ws.Range('A22').Value = 'Greene King'
$ws.Range('B22').Value = 'https://uk.indeed.com/viewjob?jk=b9bb2aaf16474db3'
$ws.Range('C22').Value = 'Posted'
$ws.Range('D22').Value = 'Competitive salary with generous bonus and pension contribution.'
$ws.Range('E22').Value = 44842
$ws.Range('E22').NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range('F22').Value = 25160
$ws.Range('G22').Value = 27676
$ws.Range('H22').Value = 'Stirling'
$ws.Range('I22').ClearContents()
